$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 5.876600000000004
$ws.Range("E3").Value = 16.2498
$ws.Range("A12").Value = -21.5648
$ws.Range("B14").Value = 6.714600000000001
$ws.Range("E20").Value = 15.92119999999999
$ws.Range("E25").Value = 17.28590000000001
$ws.Range("B26").Value = 3.867800000000004
$ws.Range("A27").Value = -21.8293
$ws.Range("E30").Value = 15.71290000000001
$ws.Range("B31").Value = 5.508000000000004
$ws.Range("A32").Value = -21.1889
$ws.Range("B35").Value = 9.462300000000004
$ws.Range("A36").Value = -19.9835
$ws.Range("B37").Value = 8.946400000000004
$ws.Range("A38").Value = -19.34070000000001
$ws.Range("E44").Value = 16.73759999999999
$ws.Range("B45").Value = 6.7126
$ws.Range("A46").Value = -21.60409999999999
$ws.Range("E47").Value = 16.43019999999999
$ws.Range("B52").Value = 5.0711
$ws.Range("A54").Value = -21.66819999999999
$ws.Range("A55").Value = -22.63590000000001
$ws.Range("A56").Value = -22.24520000000001
$ws.Range("B57").Value = 4.725599999999997
$ws.Range("E58").Value = 16.46250000000002
$ws.Range("A67").Value = -21.44739999999998
$ws.Range("A69").Value = -21.55489999999998
$ws.Range("A72").Value = -22.15720000000002
$ws.Range("E78").Value = 16.58610000000002
$ws.Range("B81").Value = 6.535800000000002
$ws.Range("A83").Value = -21.77089999999999
$ws.Range("B83").Value = 5.548300000000004
$ws.Range("E84").Value = 16.65320000000001
$ws.Range("A86").Value = -22.09990000000001
$ws.Range("E89").Value = 17.35480000000002
$ws.Range("A91").Value = -21.39190000000001
$ws.Range("E91").Value = 17.99500000000002
$ws.Range("E92").Value = 18.00340000000001
$ws.Range("A93").Value = -21.2391
$ws.Range("E96").Value = 16.18169999999999
$ws.Range("A99").Value = -20.14309999999999
$ws.Range("B100").Value = 5.224299999999996
$ws.Range("B102").Value = 8.193400000000002
$ws.Range("E102").Value = 16.73759999999998
